# Actualización automática del tracker
# - Rellena resultado/profit de las filas 113 y 114 (partidos ya resueltos)
# - Añade dos nuevos partidos (filas 118 y 119) al tracker

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Completar resultados pendientes ---
$ws.Range("G113").Value = "Acierto"
$ws.Range("H113").Value = 3

$ws.Range("G114").Value = "Fallo"
$ws.Range("H114").Value = -1

# --- Nuevas filas al final del tracker ---
# (la columna "fecha" guarda las fechas como texto "AAAA-MM-DD", no como
#  fecha-número real; forzamos formato texto y lo devolvemos a "Normal"
#  para no dejar un estilo residual en la celda)
$ws.Range("A118").Value = 14601401
$ws.Range("B118").NumberFormat = "@"
$ws.Range("B118").Value = "2025-09-14"
$ws.Range("B118").Style = "Normal"
$ws.Range("C118").Value = "Stefano Napolitano"
$ws.Range("D118").Value = "Kilian Feldbausch"
$ws.Range("E118").Value = "Gana Stefano Napolitano"
$ws.Range("F118").Value = 1.8

$ws.Range("A119").Value = 14601354
$ws.Range("B119").NumberFormat = "@"
$ws.Range("B119").Value = "2025-09-14"
$ws.Range("B119").Style = "Normal"
$ws.Range("C119").Value = "Alejandro Tabilo"
$ws.Range("D119").Value = "Juan Manuel Cerundolo"
$ws.Range("E119").Value = "Gana Juan Manuel Cerundolo"
$ws.Range("F119").Value = 2.63
